$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 1369.75
$ws.Cells.Item(6, 9).Value = 1900.4286
$ws.Cells.Item(6, 11).Value = 5701.2858
$ws.Cells.Item(6, 13).Value = -5589.2858
$ws.Cells.Item(62, 8).Value = 3645.4546
$ws.Cells.Item(62, 9).Value = 3324.5
$ws.Cells.Item(62, 11).Value = 3324.5
$ws.Cells.Item(62, 13).Value = -2700.5
$ws.Cells.Item(65, 8).Value = 3645.4546
$ws.Cells.Item(65, 9).Value = 3324.5
$ws.Cells.Item(65, 11).Value = 16622.5
$ws.Cells.Item(65, 13).Value = -13502.5
$ws.Cells.Item(111, 8).Value = 3330
$ws.Cells.Item(111, 9).Value = 2928.8
$ws.Cells.Item(111, 10).Value = 3831.5
$ws.Cells.Item(111, 11).Value = 8786.400000000001
$ws.Cells.Item(111, 12).Value = 11494.5
$ws.Cells.Item(111, 13).Value = -5719.400000000001
$ws.Cells.Item(111, 14).Value = -17628.5
$ws.Cells.Item(137, 8).Value = 50821.41
$ws.Cells.Item(137, 9).Value = 61461.715
$ws.Cells.Item(137, 11).Value = 184385.145
$ws.Cells.Item(137, 13).Value = -181835.145

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 4378.9
$ws.Cells.Item(45, 9).Value = 4612.7144
$ws.Cells.Item(45, 10).Value = 3833.3333
$ws.Cells.Item(45, 11).Value = 4612.7144
$ws.Cells.Item(45, 12).Value = 3833.3333
$ws.Cells.Item(45, 13).Value = -4235.7144
$ws.Cells.Item(45, 14).Value = -4587.3333
$ws.Cells.Item(74, 8).Value = 3274.5652
$ws.Cells.Item(74, 9).Value = 3455
$ws.Cells.Item(74, 11).Value = 3455
$ws.Cells.Item(74, 13).Value = -2581
$ws.Cells.Item(77, 8).Value = 3274.5652
$ws.Cells.Item(77, 9).Value = 3455
$ws.Cells.Item(77, 11).Value = 17275
$ws.Cells.Item(77, 13).Value = -12907
$ws.Cells.Item(102, 8).Value = 2420.077
$ws.Cells.Item(102, 9).Value = 1829
$ws.Cells.Item(102, 11).Value = 1829
$ws.Cells.Item(102, 13).Value = -207
$ws.Cells.Item(132, 8).Value = 3158.7354
$ws.Cells.Item(132, 9).Value = 2999.7666
$ws.Cells.Item(132, 11).Value = 8999.299800000001
$ws.Cells.Item(132, 13).Value = -6469.299800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2515.2856
$ws.Cells.Item(99, 9).Value = 1801.5555
$ws.Cells.Item(99, 11).Value = 1801.5555
$ws.Cells.Item(99, 13).Value = -303.5554999999999
$ws.Cells.Item(107, 8).Value = 2066.5789
$ws.Cells.Item(107, 9).Value = 670.4545000000001
$ws.Cells.Item(107, 11).Value = 670.4545000000001
$ws.Cells.Item(107, 13).Value = 1249.5455
$ws.Cells.Item(134, 8).Value = 2566342.8
$ws.Cells.Item(134, 9).Value = 2900370
$ws.Cells.Item(134, 11).Value = 8701110
$ws.Cells.Item(134, 13).Value = -8698575

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6810.8965
$ws.Cells.Item(31, 9).Value = 2850.8572
$ws.Cells.Item(31, 11).Value = 2850.8572
$ws.Cells.Item(31, 13).Value = -2555.8572
$ws.Cells.Item(34, 8).Value = 6810.8965
$ws.Cells.Item(34, 9).Value = 2850.8572
$ws.Cells.Item(34, 11).Value = 2850.8572
$ws.Cells.Item(34, 13).Value = -2648.8572
$ws.Cells.Item(58, 8).Value = 3309.9744
$ws.Cells.Item(58, 9).Value = 3093.6365
$ws.Cells.Item(58, 10).Value = 4499.8335
$ws.Cells.Item(58, 11).Value = 3093.6365
$ws.Cells.Item(58, 12).Value = 4499.8335
$ws.Cells.Item(58, 13).Value = -2890.6365
$ws.Cells.Item(58, 14).Value = -4905.8335
$ws.Cells.Item(86, 8).Value = 42017.41
$ws.Cells.Item(86, 9).Value = 33164.832
$ws.Cells.Item(86, 10).Value = 43626.97
$ws.Cells.Item(86, 11).Value = 33164.832
$ws.Cells.Item(86, 12).Value = 43626.97
$ws.Cells.Item(86, 13).Value = -32041.832
$ws.Cells.Item(86, 14).Value = -45872.97
$ws.Cells.Item(89, 8).Value = 42017.41
$ws.Cells.Item(89, 9).Value = 33164.832
$ws.Cells.Item(89, 10).Value = 43626.97
$ws.Cells.Item(89, 11).Value = 165824.16
$ws.Cells.Item(89, 12).Value = 218134.85
$ws.Cells.Item(89, 13).Value = -160208.16
$ws.Cells.Item(89, 14).Value = -229366.85
$ws.Cells.Item(94, 8).Value = 4904.8
$ws.Cells.Item(94, 9).Value = 10012
$ws.Cells.Item(94, 11).Value = 10012
$ws.Cells.Item(94, 13).Value = -9561
$ws.Cells.Item(105, 8).Value = 3374.7778
$ws.Cells.Item(105, 9).Value = 3277.6
$ws.Cells.Item(105, 10).Value = 3496.25
$ws.Cells.Item(105, 11).Value = 3277.6
$ws.Cells.Item(105, 12).Value = 3496.25
$ws.Cells.Item(105, 13).Value = -1530.6
$ws.Cells.Item(105, 14).Value = -6990.25
$ws.Cells.Item(107, 8).Value = 72525.28999999999
$ws.Cells.Item(107, 9).Value = 143513.42
$ws.Cells.Item(107, 10).Value = 1537.1428
$ws.Cells.Item(107, 11).Value = 143513.42
$ws.Cells.Item(107, 12).Value = 1537.1428
$ws.Cells.Item(107, 13).Value = -141593.42
$ws.Cells.Item(107, 14).Value = -5377.1428
$ws.Cells.Item(132, 8).Value = 4073.7585
$ws.Cells.Item(132, 9).Value = 3653.32
$ws.Cells.Item(132, 11).Value = 10959.96
$ws.Cells.Item(132, 13).Value = -8429.960000000001
$ws.Cells.Item(136, 8).Value = 3309.9744
$ws.Cells.Item(136, 9).Value = 3093.6365
$ws.Cells.Item(136, 10).Value = 4499.8335
$ws.Cells.Item(136, 11).Value = 9280.9095
$ws.Cells.Item(136, 12).Value = 13499.5005
$ws.Cells.Item(136, 13).Value = -6730.9095
$ws.Cells.Item(136, 14).Value = -18599.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 6924001
$ws.Cells.Item(121, 9).Value = 947.75
$ws.Cells.Item(121, 10).Value = 10000913
$ws.Cells.Item(121, 11).Value = 2843.25
$ws.Cells.Item(121, 12).Value = 30002739
$ws.Cells.Item(121, 13).Value = -1533.25
$ws.Cells.Item(121, 14).Value = -30005359

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 53583.332
$ws.Cells.Item(70, 9).Value = 299000
$ws.Cells.Item(70, 10).Value = 4500
$ws.Cells.Item(70, 11).Value = 299000
$ws.Cells.Item(70, 12).Value = 4500
$ws.Cells.Item(70, 13).Value = -298730
$ws.Cells.Item(70, 14).Value = -5040
$ws.Cells.Item(73, 8).Value = 53583.332
$ws.Cells.Item(73, 9).Value = 299000
$ws.Cells.Item(73, 10).Value = 4500
$ws.Cells.Item(73, 11).Value = 299000
$ws.Cells.Item(73, 12).Value = 4500
$ws.Cells.Item(73, 13).Value = -298064
$ws.Cells.Item(73, 14).Value = -6372
$ws.Cells.Item(80, 8).Value = 3200
$ws.Cells.Item(80, 9).Value = 3000
$ws.Cells.Item(80, 10).Value = 3700
$ws.Cells.Item(80, 11).Value = 3000
$ws.Cells.Item(80, 12).Value = 3700
$ws.Cells.Item(80, 13).Value = -2002
$ws.Cells.Item(80, 14).Value = -5696
$ws.Cells.Item(83, 8).Value = 3200
$ws.Cells.Item(83, 9).Value = 3000
$ws.Cells.Item(83, 10).Value = 3700
$ws.Cells.Item(83, 11).Value = 15000
$ws.Cells.Item(83, 12).Value = 18500
$ws.Cells.Item(83, 13).Value = -10008
$ws.Cells.Item(83, 14).Value = -28484
$ws.Cells.Item(122, 8).Value = 1588.6
$ws.Cells.Item(122, 9).Value = 1546.5
$ws.Cells.Item(122, 11).Value = 4639.5
$ws.Cells.Item(122, 13).Value = -2189.5
$ws.Cells.Item(126, 8).Value = 3702.625
$ws.Cells.Item(126, 10).Value = 3970.6667
$ws.Cells.Item(126, 12).Value = 11912.0001
$ws.Cells.Item(126, 14).Value = -16852.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 8405.154
$ws.Cells.Item(136, 9).Value = 4748.75
$ws.Cells.Item(136, 11).Value = 14246.25
$ws.Cells.Item(136, 13).Value = -11696.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 1000000000
$ws.Cells.Item(14, 9).Value = 1000000000
$ws.Cells.Item(14, 11).Value = 1000000000
$ws.Cells.Item(14, 13).Value = -999999832
$ws.Cells.Item(136, 8).Value = 13335215
$ws.Cells.Item(136, 9).Value = 1842.3
$ws.Cells.Item(136, 11).Value = 5526.9
$ws.Cells.Item(136, 13).Value = -2976.9
